# Updated the DegreePlan sheet

$wb = $excel.ActiveWorkbook

# --- Add the new DegreePlan rows (for students S521315 and S533985) ---
$ws = $wb.Worksheets.Item("DegreePlan")
$ws.Activate()

$ws.Range("A4").Value = 10
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "S521315"
$ws.Range("D4").Value = "No summer off"
$ws.Range("E4").Value = "No summer off"

$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "S521315"
$ws.Range("D5").Value = "Summer Off"
$ws.Range("E5").Value = "summer off"

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "S533985"
$ws.Range("D6").Value = "No summer off"
$ws.Range("E6").Value = "No summer off"

$ws.Range("A7").Value = 11
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "S533985"
$ws.Range("D7").Value = "Summer Off"
$ws.Range("E7").Value = "summer off"

# Selection on the DegreePlan sheet ends up at B8
$ws.Range("B8").Select()

$wb.Save()
